$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 209
$ws.Range("J6").Value = 37
$ws.Range("L6").Value = 111
$ws.Range("N6").Value = -335
$ws.Range("H31").Value = 56.666668
$ws.Range("I31").Value = 56.666668
$ws.Range("K31").Value = 170.000004
$ws.Range("M31").Value = 59.99999600000001
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").ClearContents()
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = 0
$ws.Range("H51").Value = 6656.3335
$ws.Range("I51").Value = 6499.5
$ws.Range("J51").Value = 6970
$ws.Range("K51").Value = 6499.5
$ws.Range("L51").Value = 6970
$ws.Range("M51").Value = -6015.5
$ws.Range("N51").Value = -7938
$ws.Range("H62").Value = 4383.6665
$ws.Range("I62").Value = 4383.6665
$ws.Range("K62").Value = 4383.6665
$ws.Range("M62").Value = -3759.6665
$ws.Range("H65").Value = 4383.6665
$ws.Range("I65").Value = 4383.6665
$ws.Range("K65").Value = 21918.3325
$ws.Range("M65").Value = -18798.3325
$ws.Range("H106").Value = 2839.8
$ws.Range("I106").Value = 2839.8
$ws.Range("K106").Value = 2839.8
$ws.Range("M106").Value = -2208.8
$ws.Range("H113").Value = 9154.777
$ws.Range("I113").Value = 9469.6
$ws.Range("K113").Value = 9469.6
$ws.Range("M113").Value = -6215.6
$ws.Range("H132").Value = 157476.08
$ws.Range("I132").Value = 4025.25
$ws.Range("J132").Value = 402997.4
$ws.Range("K132").Value = 12075.75
$ws.Range("L132").Value = 1208992.2
$ws.Range("M132").Value = -9545.75
$ws.Range("N132").Value = -1214052.2
$ws.Range("H137").Value = 3409.1177
$ws.Range("I137").Value = 2117.1667
$ws.Range("K137").Value = 6351.500100000001
$ws.Range("M137").Value = -3801.500100000001
$ws.Range("H138").Value = 8586.3125
$ws.Range("J138").Value = 9132.75
$ws.Range("L138").Value = 27398.25
$ws.Range("N138").Value = -37678.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2969.2
$ws.Range("I20").Value = 2969.2
$ws.Range("K20").Value = 2969.2
$ws.Range("M20").Value = -2722.2
$ws.Range("H58").Value = 88419.875
$ws.Range("J58").Value = 92481.86
$ws.Range("L58").Value = 92481.86
$ws.Range("N58").Value = -93069.86
$ws.Range("H107").Value = 2922.6924
$ws.Range("I107").Value = 2869.625
$ws.Range("J107").Value = 3007.6
$ws.Range("K107").Value = 2869.625
$ws.Range("L107").Value = 3007.6
$ws.Range("M107").Value = -949.625
$ws.Range("N107").Value = -6847.6

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4389.087
$ws.Range("I31").Value = 1446.4
$ws.Range("J31").Value = 5206.5
$ws.Range("K31").Value = 1446.4
$ws.Range("L31").Value = 5206.5
$ws.Range("M31").Value = -1151.4
$ws.Range("N31").Value = -5796.5
$ws.Range("H34").Value = 4389.087
$ws.Range("I34").Value = 1446.4
$ws.Range("J34").Value = 5206.5
$ws.Range("K34").Value = 1446.4
$ws.Range("L34").Value = 5206.5
$ws.Range("M34").Value = -1244.4
$ws.Range("N34").Value = -5610.5
$ws.Range("H58").Value = 5977.1
$ws.Range("I58").Value = 5902.8
$ws.Range("J58").Value = 6200
$ws.Range("K58").Value = 5902.8
$ws.Range("L58").Value = 6200
$ws.Range("M58").Value = -5699.8
$ws.Range("N58").Value = -6606
$ws.Range("H134").Value = 3718.2
$ws.Range("I134").Value = 3776.7856
$ws.Range("K134").Value = 11330.3568
$ws.Range("M134").Value = -8795.356800000001
$ws.Range("H136").Value = 5977.1
$ws.Range("I136").Value = 5902.8
$ws.Range("J136").Value = 6200
$ws.Range("K136").Value = 17708.4
$ws.Range("L136").Value = 18600
$ws.Range("M136").Value = -15158.4
$ws.Range("N136").Value = -23700

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 519.5
$ws.Range("I17").Value = 40
$ws.Range("J17").Value = 999
$ws.Range("K17").Value = 120
$ws.Range("L17").Value = 2997
$ws.Range("M17").Value = 49
$ws.Range("N17").Value = -3335
$ws.Range("H58").Value = 500
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H68").Value = 2921.682
$ws.Range("I68").Value = 1300
$ws.Range("J68").Value = 2998.9048
$ws.Range("K68").Value = 3900
$ws.Range("L68").Value = 8996.714399999999
$ws.Range("M68").Value = -3089
$ws.Range("N68").Value = -10618.7144
$ws.Range("H71").Value = 2921.682
$ws.Range("I71").Value = 1300
$ws.Range("J71").Value = 2998.9048
$ws.Range("K71").Value = 11700
$ws.Range("L71").Value = 26990.1432
$ws.Range("M71").Value = -7644
$ws.Range("N71").Value = -35102.1432
$ws.Range("H80").Value = 10099.889
$ws.Range("I80").Value = 10239.8
$ws.Range("J80").Value = 9925
$ws.Range("K80").Value = 30719.4
$ws.Range("L80").Value = 29775
$ws.Range("M80").Value = -29783.4
$ws.Range("N80").Value = -31647
$ws.Range("H83").Value = 10099.889
$ws.Range("I83").Value = 10239.8
$ws.Range("J83").Value = 9925
$ws.Range("K83").Value = 92158.2
$ws.Range("L83").Value = 89325
$ws.Range("M83").Value = -87478.2
$ws.Range("N83").Value = -98685
$ws.Range("H107").Value = 1573.8182
$ws.Range("I107").Value = 1256.5714
$ws.Range("J107").Value = 1721.8667
$ws.Range("K107").Value = 3769.7142
$ws.Range("L107").Value = 5165.6001
$ws.Range("M107").Value = -1849.7142
$ws.Range("N107").Value = -9005.6001
$ws.Range("H131").Value = 20099.383
$ws.Range("J131").Value = 1655.6182
$ws.Range("L131").Value = 4966.8546
$ws.Range("N131").Value = -15046.8546

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 332.7
$ws.Range("I2").Value = 258.55554
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 258.55554
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = -145.55554
$ws.Range("N2").Value = -1226
$ws.Range("H70").Value = 8525
$ws.Range("I70").Value = 7258.9
$ws.Range("J70").Value = 9931.777
$ws.Range("K70").Value = 7258.9
$ws.Range("L70").Value = 9931.777
$ws.Range("M70").Value = -6988.9
$ws.Range("N70").Value = -10471.777
$ws.Range("H73").Value = 8525
$ws.Range("I73").Value = 7258.9
$ws.Range("J73").Value = 9931.777
$ws.Range("K73").Value = 7258.9
$ws.Range("L73").Value = 9931.777
$ws.Range("M73").Value = -6322.9
$ws.Range("N73").Value = -11803.777
$ws.Range("H102").Value = 4125.5
$ws.Range("J102").Value = 3796
$ws.Range("L102").Value = 3796
$ws.Range("N102").Value = -7040
$ws.Range("H113").Value = 5000
$ws.Range("I113").Value = 5000
$ws.Range("K113").Value = 5000
$ws.Range("M113").Value = -2830
$ws.Range("H126").Value = 3860
$ws.Range("I126").Value = 2389
$ws.Range("J126").Value = 4350.3335
$ws.Range("K126").Value = 7167
$ws.Range("L126").Value = 13051.0005
$ws.Range("M126").Value = -4697
$ws.Range("N126").Value = -17991.0005
$ws.Range("H136").Value = 67130.25
$ws.Range("J136").Value = 67130.25
$ws.Range("L136").Value = 201390.75
$ws.Range("N136").Value = -206490.75

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3756.1052
$ws.Range("I132").Value = 2540.5
$ws.Range("K132").Value = 7621.5
$ws.Range("M132").Value = -5091.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4277.4
$ws.Range("I81").Value = 3462.3333
$ws.Range("J81").Value = 5500
$ws.Range("K81").Value = 6924.6666
$ws.Range("L81").Value = 11000
$ws.Range("M81").Value = -5863.6666
$ws.Range("N81").Value = -13122
$ws.Range("H84").Value = 4277.4
$ws.Range("I84").Value = 3462.3333
$ws.Range("J84").Value = 5500
$ws.Range("K84").Value = 34623.333
$ws.Range("L84").Value = 55000
$ws.Range("M84").Value = -29319.333
$ws.Range("N84").Value = -65608
